$wb = $excel.ActiveWorkbook

# Rename sheets to reflect the new "strain_log2_expression" naming scheme
$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Remove tabSelected from the previously-active sheet (optimization_parameters)
# by activating the now-current sheet instead.
$wsDcin5.Activate() | Out-Null
$wsDcin5.Range("F43").Select() | Out-Null
